# "Generate Report for Archive"
#
# The localization-status report is regenerated: the Status of the single
# tracked file flips from "Ready for handoff" to "In Translation" on every
# sheet that surfaces it (the zh-cn/de-de per-language tables' "Status"
# column, and the corresponding zh-cn/de-de columns on the Overview roll-up
# sheet). Excel then re-autosizes the now-narrower Status columns to fit the
# new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 13.4101845877511

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
